$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C21").Value = "dsgagass"
$ws.Range("C22").Value = "dsvsbsb"
